# "Cap nhat ten bai tap" - update exercise name / title text across the deck.
$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 1: merge "Bai 14" + ". Phan lop van ban" into a single run.
# ---------------------------------------------------------------------------
$s1 = $p.Slides.Item(1)
$shTitle = $s1.Shapes.Item(2)
$trTitle = $shTitle.TextFrame.TextRange
$para1 = $trTitle.Paragraphs(1,1)
# Force a full rewrite of the run by first replacing with a throwaway string,
# then setting the real text - this collapses the paragraph back to one run.
$para1.Text = "TEMP_PLACEHOLDER"
$trTitle2 = $shTitle.TextFrame.TextRange
$para1b = $trTitle2.Paragraphs(1,1)
$para1b.Text = "Bài 14. Phân lớp văn bản"

# ---------------------------------------------------------------------------
# Slide 28: "Bai tap (2)" -> "Bai tap 14.1", resize content box, and split the
# paragraph in two.
# ---------------------------------------------------------------------------
$s28 = $p.Slides.Item(28)

# -- Title shape: replace the trailing " (2)" run with " " + "14.1" runs.
$shT = $s28.Shapes.Item(1)
$trT = $shT.TextFrame.TextRange
$tail = $trT.Characters(8,4)
$tail.Text = "ZZZZZZZZ"
$trT2 = $shT.TextFrame.TextRange
$tail2 = $trT2.Characters(8,8)
$tail2.Text = " 14.1"
$trT3 = $shT.TextFrame.TextRange
$numPart = $trT3.Characters(9,4)
$numPart.Text = "YYYY"
$trT4 = $shT.TextFrame.TextRange
$numPart2 = $trT4.Characters(9,4)
$numPart2.Text = "14.1"

# -- Content placeholder shape: reposition/resize and split the paragraph.
$shC = $s28.Shapes.Item(2)
$shC.Left = 59.494173228346455
$shC.Top = 158.87503937007875
$shC.Width = 645.6309211417307
$shC.Height = 324.0

$trC = $shC.TextFrame.TextRange
$fullText = $trC.Text
$splitAt = 163
$newText = $fullText.Substring(0, $splitAt) + "`r" + $fullText.Substring($splitAt)
$trC.Text = $newText

# Re-split "Hay " into its own run at the start of the new second paragraph.
$trC2 = $shC.TextFrame.TextRange
$hayTemp = $trC2.Characters($splitAt + 2, 4)
$hayTemp.Text = "ZZZZZZZZ"
$trC3 = $shC.TextFrame.TextRange
$hayFinal = $trC3.Characters($splitAt + 2, 8)
$hayFinal.Text = "Hãy "
